$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Qty/Unit-cost values for rows 11-15
$ws.Range("G11").Value = 0.2
$ws.Range("G12").Value = 0.01
$ws.Range("G13").Value = 0.2
$ws.Range("G14").Value = 1.6
$ws.Range("G15").Value = 0.14000000000000001

# Total cost ($) formulas for rows 5-15: (Qty * Unit price)
$ws.Range("H5").Formula = "=(B5*G5)"
$ws.Range("H6:H15").Formula = "=(B6*G6)"

# Grand total row
$ws.Range("G16").Value = "GRAND TOTAL PER HAT ($)"
$ws.Range("G16").Font.Bold = $true
$ws.Range("H16").Formula = "=SUM(H5:H15)"

$null = $ws.Range("I16").Select()

$ws.PageSetup.Orientation = 1
